# The document has two logo pictures that each appear twice (once in the
# "primary" header/footer, once in the "first page" header/footer):
#   - BTec_Logo-Orange        -> lives in both headers
#   - PearsonLogo.png (descr) -> lives in both footers
#
# The edit simply renames the inline picture (its <wp:docPr name="...">,
# i.e. the Word "Name" shown in Selection Pane / Format tab):
#   headers: image1.jpg -> image2.jpg
#   footers: image2.png -> image1.png
#
# Renaming through Range.InlineShapes(1).Name directly throws
# "addressed block not found" for shapes that live in a footer story, so we
# first Select() the shape and then rename it through $word.Selection -
# that re-resolves the shape against a fresh (non-stale) anchor and works
# for both headers and footers.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlineLogo($headerOrFooter, [string]$newName) {
    if ($headerOrFooter.Exists -and $headerOrFooter.Range.InlineShapes.Count -ge 1) {
        $shp = $headerOrFooter.Range.InlineShapes(1)
        $shp.Select() | Out-Null
        $word.Selection.InlineShapes(1).Name = $newName
    }
}

# Headers: BTec_Logo-Orange, image1.jpg -> image2.jpg
for ($i = 1; $i -le 3; $i++) {
    Rename-InlineLogo $sec.Headers($i) "image2.jpg"
}

# Footers: PearsonLogo.png, image2.png -> image1.png
for ($i = 1; $i -le 3; $i++) {
    Rename-InlineLogo $sec.Footers($i) "image1.png"
}
